# Update "想去人数" (want-to-go count) figures to the freshly scraped
# values (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet -------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F8").Value  = 10926
$wsExpo.Range("F9").Value  = 193
$wsExpo.Range("F19").Value = 1165
$wsExpo.Range("F27").Value = 3297
$wsExpo.Range("F39").Value = 1342
$wsExpo.Range("F40").Value = 2320
$wsExpo.Range("F41").Value = 5411
$wsExpo.Range("F47").Value = 15

# --- 演出 (Performances) sheet ------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value  = 4116
$wsShow.Range("F11").Value = 462

# --- 本地生活 (Local life) sheet -----------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 52

# --- 全部类型 (All types) sheet ------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 52
$wsAll.Range("F8").Value  = 10926
$wsAll.Range("F9").Value  = 193
$wsAll.Range("F15").Value = 1165
$wsAll.Range("F18").Value = 4116
$wsAll.Range("F24").Value = 3297
$wsAll.Range("F34").Value = 2322
$wsAll.Range("F42").Value = 15
